$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header relabel: average_doctor / average_doctor_old columns swapped (BP <-> BQ)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Harvard case classification recalculated statistics (rows 4-13)
# Row 4
$ws.Range("AI4").Value = 0.286
$ws.Range("AJ4").Value = 0.089
$ws.Range("AK4").Value = 0.298
$ws.Range("AU4").Value = 0.187
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.168
$ws.Range("BA4").Value = 2.005
$ws.Range("BB4").Value = 0.158
$ws.Range("BC4").Value = 0.397
$ws.Range("BG4").Value = 0.729
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.376
$ws.Range("BM4").Value = 0.72
$ws.Range("BN4").Value = 0.081
$ws.Range("BO4").Value = 0.284
$ws.Range("BP4").Value = 0.668
$ws.Range("BQ4").Value = 0.709
$ws.Range("E4").Value = 0.424
$ws.Range("F4").Value = 0.072
$ws.Range("G4").Value = 0.269
$ws.Range("N4").Value = 0.434
$ws.Range("O4").Value = 0.064
$ws.Range("P4").Value = 0.254
$ws.Range("Q4").Value = 0.024
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.129
$ws.Range("W4").Value = 0.286
$ws.Range("X4").Value = 0.109
$ws.Range("Y4").Value = 0.33
# Row 5
$ws.Range("AI5").Value = 0.309
$ws.Range("AJ5").Value = 0.1
$ws.Range("AK5").Value = 0.316
$ws.Range("AU5").Value = 0.366
$ws.Range("AV5").Value = 0.094
$ws.Range("AW5").Value = 0.307
$ws.Range("BA5").Value = 1.337
$ws.Range("BC5").Value = 0.282
$ws.Range("BG5").Value = 0.399
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.227
$ws.Range("BM5").Value = 0.547
$ws.Range("BN5").Value = 0.064
$ws.Range("BO5").Value = 0.253
$ws.Range("BP5").Value = 0.446
$ws.Range("BQ5").Value = 0.456
$ws.Range("E5").Value = 0.541
$ws.Range("F5").Value = 0.086
$ws.Range("G5").Value = 0.293
$ws.Range("N5").Value = 0.737
$ws.Range("O5").Value = 0.08
$ws.Range("P5").Value = 0.283
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.083
$ws.Range("W5").Value = 0.276
$ws.Range("X5").Value = 0.109
$ws.Range("Y5").Value = 0.33
# Row 6
$ws.Range("AI6").Value = 0.297
$ws.Range("AU6").Value = 0.248
$ws.Range("BA6").Value = 1.597
$ws.Range("BG6").Value = 0.516
$ws.Range("BM6").Value = 0.622
$ws.Range("BP6").Value = 0.532
$ws.Range("BQ6").Value = 0.552
$ws.Range("E6").Value = 0.475
$ws.Range("N6").Value = 0.546
$ws.Range("Q6").Value = 0.019
$ws.Range("W6").Value = 0.281
# Row 7
$ws.Range("AI7").Value = 0.304
$ws.Range("AU7").Value = 0.307
$ws.Range("BA7").Value = 1.43
$ws.Range("BG7").Value = 0.439
$ws.Range("BM7").Value = 0.575
$ws.Range("BP7").Value = 0.477
$ws.Range("BQ7").Value = 0.49
$ws.Range("E7").Value = 0.513
$ws.Range("N7").Value = 0.647
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.278
# Row 8
$ws.Range("AI8").Value = 0.332
$ws.Range("AJ8").Value = 0.13
$ws.Range("AK8").Value = 0.361
$ws.Range("AU8").Value = 0.309
$ws.Range("AV8").Value = 0.085
$ws.Range("AW8").Value = 0.292
$ws.Range("BA8").Value = 1.753
$ws.Range("BB8").Value = 0.125
$ws.Range("BC8").Value = 0.354
$ws.Range("BG8").Value = 0.57
$ws.Range("BH8").Value = 0.107
$ws.Range("BI8").Value = 0.326
$ws.Range("BM8").Value = 0.69
$ws.Range("BN8").Value = 0.067
$ws.Range("BO8").Value = 0.259
$ws.Range("BP8").Value = 0.584
$ws.Range("BQ8").Value = 0.606
$ws.Range("E8").Value = 0.603
$ws.Range("F8").Value = 0.113
$ws.Range("G8").Value = 0.335
$ws.Range("N8").Value = 0.779
$ws.Range("O8").Value = 0.063
$ws.Range("P8").Value = 0.251
$ws.Range("Q8").Value = 0.017
$ws.Range("W8").Value = 0.304
# Row 9
$ws.Range("AI9").Value = 0.258
$ws.Range("AJ9").Value = 0.191
$ws.Range("AK9").Value = 0.438
$ws.Range("BA9").Value = 1.731
$ws.Range("BB9").Value = 0.249
$ws.Range("BC9").Value = 0.499
$ws.Range("BG9").Value = 0.613
$ws.Range("BH9").Value = 0.237
$ws.Range("BI9").Value = 0.487
$ws.Range("BM9").Value = 0.656
$ws.Range("BN9").Value = 0.226
$ws.Range("BO9").Value = 0.475
$ws.Range("BP9").Value = 0.577
$ws.Range("BQ9").Value = 0.593
$ws.Range("E9").Value = 0.548
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.688
$ws.Range("O9").Value = 0.215
$ws.Range("P9").Value = 0.463
$ws.Range("W9").Value = 0.204
$ws.Range("X9").Value = 0.163
$ws.Range("Y9").Value = 0.403
# Row 10
$ws.Range("AI10").Value = 0.366
$ws.Range("AJ10").Value = 0.232
$ws.Range("AK10").Value = 0.482
$ws.Range("AU10").Value = 0.301
$ws.Range("AV10").Value = 0.21
$ws.Range("AW10").Value = 0.459
$ws.Range("BA10").Value = 2.097
$ws.Range("BB10").Value = 0.242
$ws.Range("BC10").Value = 0.492
$ws.Range("BG10").Value = 0.667
$ws.Range("BH10").Value = 0.222
$ws.Range("BI10").Value = 0.471
$ws.Range("BM10").Value = 0.839
$ws.Range("BN10").Value = 0.135
$ws.Range("BO10").Value = 0.368
$ws.Range("BP10").Value = 0.699
$ws.Range("BQ10").Value = 0.729
$ws.Range("E10").Value = 0.677
$ws.Range("F10").Value = 0.219
$ws.Range("G10").Value = 0.467
$ws.Range("N10").Value = 0.882
$ws.Range("O10").Value = 0.104
$ws.Range("P10").Value = 0.323
$ws.Range("W10").Value = 0.376
$ws.Range("X10").Value = 0.235
$ws.Range("Y10").Value = 0.484
# Row 11
$ws.Range("AI11").Value = 0.398
$ws.Range("AJ11").Value = 0.24
$ws.Range("AK11").Value = 0.489
$ws.Range("AU11").Value = 0.441
$ws.Range("AV11").Value = 0.247
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.097
$ws.Range("BB11").Value = 0.242
$ws.Range("BC11").Value = 0.492
$ws.Range("BG11").Value = 0.667
$ws.Range("BH11").Value = 0.222
$ws.Range("BI11").Value = 0.471
$ws.Range("BM11").Value = 0.839
$ws.Range("BN11").Value = 0.135
$ws.Range("BO11").Value = 0.368
$ws.Range("BP11").Value = 0.699
$ws.Range("BQ11").Value = 0.732
$ws.Range("E11").Value = 0.71
$ws.Range("F11").Value = 0.206
$ws.Range("G11").Value = 0.454
$ws.Range("N11").Value = 0.903
$ws.Range("O11").Value = 0.087
$ws.Range("P11").Value = 0.296
$ws.Range("W11").Value = 0.376
$ws.Range("X11").Value = 0.235
$ws.Range("Y11").Value = 0.484
# Row 12
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.767
$ws.Range("AV12").Value = 2.737
$ws.Range("AW12").Value = 1.654
$ws.Range("BA12").Value = 3.688
$ws.Range("BB12").Value = 0.395
$ws.Range("BC12").Value = 0.629
$ws.Range("BG12").Value = 1.097
$ws.Range("BH12").Value = 0.12
$ws.Range("BI12").Value = 0.346
$ws.Range("BM12").Value = 1.282
$ws.Range("BN12").Value = 0.331
$ws.Range("BO12").Value = 0.575
$ws.Range("BP12").Value = 1.229
$ws.Range("BQ12").Value = 1.255
$ws.Range("E12").Value = 1.409
$ws.Range("F12").Value = 0.757
$ws.Range("G12").Value = 0.87
$ws.Range("N12").Value = 1.465
$ws.Range("O12").Value = 1.039
$ws.Range("P12").Value = 1.02
$ws.Range("W12").Value = 1.629
$ws.Range("X12").Value = 0.576
$ws.Range("Y12").Value = 0.759
# Row 13
$ws.Range("AI13").Value = 1.288
$ws.Range("AJ13").Value = 0.37
$ws.Range("AK13").Value = 0.608
$ws.Range("AU13").Value = 2.296
$ws.Range("AV13").Value = 0.925
$ws.Range("AW13").Value = 0.962
$ws.Range("BA13").Value = 2.352
$ws.Range("BB13").Value = 0.296
$ws.Range("BC13").Value = 0.544
$ws.Range("BG13").Value = 0.587
$ws.Range("BH13").Value = 0.071
$ws.Range("BI13").Value = 0.266
$ws.Range("BM13").Value = 0.891
$ws.Range("BN13").Value = 0.285
$ws.Range("BO13").Value = 0.533
$ws.Range("BP13").Value = 0.784
$ws.Range("BQ13").Value = 0.726
$ws.Range("E13").Value = 1.579
$ws.Range("F13").Value = 0.656
$ws.Range("G13").Value = 0.81
$ws.Range("N13").Value = 2.069
$ws.Range("O13").Value = 0.938
$ws.Range("P13").Value = 0.969
$ws.Range("W13").Value = 1.037
$ws.Range("X13").Value = 0.193
$ws.Range("Y13").Value = 0.439
